# Bugfixed the naive forecaster component module
#
# The "date" column (column A, rows 2-82) previously stored raw Excel
# date serials (displayed via a custom "YYYY-MM-DD HH:MM:SS" number
# format). The naive forecaster is supposed to key its series off
# calendar quarters, not literal dates, so each date is converted to a
# "<year>Q<quarter>" text label (e.g. 38398 -> "2005-02-15" -> "2005Q1").
#
# After the conversion the cells hold plain text, so their formatting is
# re-based on the header cell's style (bold font, thin border, centered
# alignment, default/General number format) instead of the old date
# format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Excel's date serial epoch (serial 0 == 1899-12-30).
$epoch = Get-Date -Year 1899 -Month 12 -Day 30

$firstRow = 2
$lastRow = $ws.UsedRange.Rows.Count

for ($row = $firstRow; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 1)
    $serial = $cell.Value2
    $dt = $epoch.AddDays($serial)
    $quarter = [Math]::Floor(($dt.Month - 1) / 3) + 1
    $cell.Value = "{0}Q{1}" -f $dt.Year, $quarter
}

# Re-apply the header's cell style (bold font, thin border, centered
# alignment, General number format) to the whole label column so it no
# longer carries the old date number format.
$ws.Range("A1").Copy() | Out-Null
$ws.Range("A" + $firstRow + ":A" + $lastRow).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
